# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet right after "总计" (pushing every later
# quarter sheet one slot back), fills it with the Q3 fund-holding table, and
# inserts a matching summary row at the top of the "总计" data table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)          # "总计"

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# NOTE: sheet handles in this host track tab *position*, not identity, so
# the "2022-Q2" lookup must happen AFTER the new tab is inserted in front of
# it (otherwise it would still be bound to the slot the new sheet just took
# over). We reuse it both for the header-row style and the column-A style.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")  # existing sheet we borrow formatting from

# Header row (identical labels/style to every other quarter sheet) - copy
# straight from 2022-Q2 so the style index matches exactly.
$q2Sheet.Range("B1:H1").Copy($q3Sheet.Range("B1:H1"))


# ---------------------------------------------------------------------------
# 2. Fill in the 16 fund rows (A=index, B..G=text fields, H=numeric rank)
# ---------------------------------------------------------------------------
$q3rows = @(
    @("512290", "国泰中证生物医药ETF", "40.70", "99.74", "2.15", "0.8750", "9"),
    @("011201", "财通优势行业轮动混合A", "8.53", "91.75", "4.15", "0.3540", "9"),
    @("161122", "易方达中证万得生物科技指数（LOF）A", "8.94", "94.64", "2.25", "0.2012", "8"),
    @("320018", "诺安新动力灵活配置混合A", "4.62", "52.38", "1.95", "0.0901", "9"),
    @("000326", "南方中小盘成长股票", "4.02", "91.10", "1.78", "0.0716", "7"),
    @("013599", "华润元大臻选回报混合C", "0.78", "84.82", "6.10", "0.0476", "8"),
    @("014551", "诺安新动力灵活配置混合C", "2.32", "52.38", "1.95", "0.0452", "9"),
    @("167506", "安信深圳科技指数（LOF）A", "0.82", "93.20", "3.70", "0.0303", "9"),
    @("005117", "金信价值精选灵活配置混合A", "0.72", "87.60", "3.58", "0.0258", "6"),
    @("010572", "易方达中证万得生物科技指数（LOF）C", "0.99", "94.64", "2.25", "0.0223", "8"),
    @("011202", "财通优势行业轮动混合C", "0.34", "91.75", "4.15", "0.0141", "9"),
    @("167507", "安信深圳科技指数（LOF）C", "0.30", "93.20", "3.70", "0.0111", "9"),
    @("005118", "金信价值精选灵活配置混合C", "0.05", "87.60", "3.58", "0.0018", "6"),
    @("519222", "海富通欣益灵活配置混合A", "0.29", "23.65", "0.32", "0.0009", "10"),
    @("519221", "海富通欣益灵活配置混合C", "0.09", "23.65", "0.32", "0.0003", "10"),
    @("013598", "华润元大臻选回报混合A", "0.00", "84.82", "6.10", "0", "8")
)

# B..G hold text (fund codes keep leading zeros, decimals keep trailing
# zeros exactly as scraped) - pre-format as Text so Excel doesn't silently
# coerce them to numbers on entry.
$q3Sheet.Range("B2:G17").NumberFormat = "@"

for ($r = 0; $r -lt $q3rows.Length; $r++) {
    $row = $q3rows[$r]
    $excelRow = $r + 2

    $q3Sheet.Cells.Item($excelRow, 1).Value = $r             # A: 0-based index
    $q3Sheet.Cells.Item($excelRow, 2).Value = $row[0]         # B: 基金代码
    $q3Sheet.Cells.Item($excelRow, 3).Value = $row[1]         # C: 基金名称
    $q3Sheet.Cells.Item($excelRow, 4).Value = $row[2]         # D: 基金规模
    $q3Sheet.Cells.Item($excelRow, 5).Value = $row[3]         # E: 股票总仓位
    $q3Sheet.Cells.Item($excelRow, 6).Value = $row[4]         # F: 仓位占比
    $q3Sheet.Cells.Item($excelRow, 7).Value = $row[5]         # G: 持有市值(亿元)
    $q3Sheet.Cells.Item($excelRow, 8).Value = [double]$row[6] # H: 仓位排名 (numeric)
}

# Row 17 (华润元大臻选回报混合A) has an exact-zero holding value, recorded as
# a plain number instead of text in the source data.
$q3Sheet.Range("G17").ClearFormats()
$q3Sheet.Range("G17").Value = 0

# Column A (0-based index) is styled the same bordered/bold style used by
# every other quarter sheet - copy it over from the neighbouring 2022-Q2 tab.
$q2Sheet.Range("A2").Copy()
$q3Sheet.Range("A2:A17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Insert the new summary row at the top of "总计"'s data table
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 1.79

# Re-apply the bordered/bold "index" style to A2 (Insert() doesn't carry it
# down), matching A3:A9 below it; then strip the stray format Insert() left
# behind on B2:D2 so they stay plain like the rest of the table.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$totalSheet.Range("B2:D2").ClearFormats()
$excel.CutCopyMode = $false
